$wb = $excel.ActiveWorkbook

# Update "F" column (想去人数) values on both the "展览" and "全部类型" sheets.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 144
    $ws.Range("F3").Value = 221
    $ws.Range("F4").Value = 3742
}
